# Update the "Analysis Results" dashboard sheet:
#  - Insert a new "Date and Time" row at the very top.
#  - Insert a new "Total SOC consumed(%)" row and a new "Cycle Count of
#    battery" row, relabel several parameters with explicit units, fix
#    the Highest/Lowest Cell Voltage and highest/lowest cell temp swap,
#    update computed metrics, and append three new speed-bucket rows
#    (70-80 km/h, 80-90 km/h).
#
# Net effect: dimension grows from A1:B42 to A1:B45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Date and Time"
$ws.Range("B1").Value = "2024-03-12 11:35:09.335000 to 2024-03-12 13:29:38.656000"

$ws.Range("A2").Value = "Total time taken for the ride"
$ws.Range("B2").Value = 0.07951575231481482
$ws.Range("B2").NumberFormat = "[hh]:mm:ss"

$ws.Range("A3").Value = "Actual Ampere-hours (Ah)"
$ws.Range("B3").Value = 33.27927277777778

$ws.Range("A4").Value = "Actual Watt-hours (Wh)"
$ws.Range("B4").Value = 1717.139885148611

$ws.Range("A5").Value = "Starting SoC (Ah)"
$ws.Range("B5").Value = 4.697

$ws.Range("A6").Value = "Ending SoC (Ah)"
$ws.Range("B6").Value = 4.135

$ws.Range("A7").Value = "Starting SoC (%)"
$ws.Range("B7").Value = 10

$ws.Range("A8").Value = "Ending SoC (%)"
$ws.Range("B8").Value = 100

$ws.Range("A9").Value = "Total distance covered (km)"
$ws.Range("B9").Value = 56.04125756068157

$ws.Range("A10").Value = "Total energy consumption(WH/KM)"
$ws.Range("B10").Value = 30.64063798513602

$ws.Range("A11").Value = "Total SOC consumed(%)"
$ws.Range("B11").Value = 90

$ws.Range("A12").Value = "Mode"
$ws.Range("B12").Value = "Eco mode`n99.75%`nSports mode`n0.04%`nCustom mode`n0.04%"
$ws.Rows.Item(12).EntireRow.AutoFit()

$ws.Range("A13").Value = "Peak Power(kW)"
$ws.Range("B13").Value = 2471.347216

$ws.Range("A14").Value = "Average Power(kW)"
$ws.Range("B14").Value = -904.6836801602516

$ws.Range("A15").Value = "Total Energy Regenerated(kWh)"
$ws.Range("B15").Value = 79.19017365777778

$ws.Range("A16").Value = "Regenerative Effectiveness(%)"
$ws.Range("B16").Value = 4.408442272039773

$ws.Range("A17").Value = "Highest Cell Voltage(V)"
$ws.Range("B17").Value = 3.382

$ws.Range("A18").Value = "Lowest Cell Voltage(V)"
$ws.Range("B18").Value = 3.125

$ws.Range("A19").Value = "Difference in Cell Voltage(V)"
$ws.Range("B19").Value = 0.2570000000000001

$ws.Range("A20").Value = "Minimum Temperature(C)"
$ws.Range("B20").Value = 28

$ws.Range("A21").Value = "Maximum Temperature(C)"
$ws.Range("B21").Value = 46

$ws.Range("A22").Value = "Difference in Temperature(C)"
$ws.Range("B22").Value = 18

$ws.Range("A23").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("B23").Value = 57

$ws.Range("A24").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("B24").Value = 61

$ws.Range("A25").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("B25").Value = 59

$ws.Range("A26").Value = "Maximum MCU Temperature(C)"
$ws.Range("B26").Value = 45

$ws.Range("A27").Value = "Maximum Motor Temperature(C)"
$ws.Range("B27").Value = 0

$ws.Range("A28").Value = "Abnormal Motor Temperature Detected(C)"
$ws.Range("B28").Value = 0

$ws.Range("A29").Value = "highest cell temp(C)"
$ws.Range("B29").Value = 46

$ws.Range("A30").Value = "lowest cell temp(C)"
$ws.Range("B30").Value = 28

$ws.Range("A31").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Range("B31").Value = 18

$ws.Range("A32").Value = "Battery Voltage(V)"
$ws.Range("B32").Value = 55

$ws.Range("A33").Value = "Total energy charged(kWh)"
$ws.Range("B33").Value = 1.830360002777778

$ws.Range("A34").Value = "Electricity consumption units(kW)"
$ws.Range("B34").Value = 0.000000074007763334052149563812

$ws.Range("A35").Value = "Cycle Count of battery"
$ws.Range("B35").Value = 115

$ws.Range("A36").Value = "Idling time percentage"
$ws.Range("B36").Value = 2.723747980613894

$ws.Range("A37").Value = "Time spent in 0-10 km/h"
$ws.Range("B37").Value = 18.17770597738288

$ws.Range("A38").Value = "Time spent in 10-20 km/h"
$ws.Range("B38").Value = 3.647819063004846

$ws.Range("A39").Value = "Time spent in 20-30 km/h"
$ws.Range("B39").Value = 7.846526655896607

$ws.Range("A40").Value = "Time spent in 30-40 km/h"
$ws.Range("B40").Value = 45.86914378029079

$ws.Range("A41").Value = "Time spent in 40-50 km/h"
$ws.Range("B41").Value = 20.29886914378029

$ws.Range("A42").Value = "Time spent in 50-60 km/h"
$ws.Range("B42").Value = 0

$ws.Range("A43").Value = "Time spent in 60-70 km/h"
$ws.Range("B43").Value = 0

$ws.Range("A44").Value = "Time spent in 70-80 km/h"
$ws.Range("B44").Value = 0

$ws.Range("A45").Value = "Time spent in 80-90 km/h"
$ws.Range("B45").Value = 0

Write-Host "Edit complete"
